$p = $ppt.ActivePresentation

# The author removed the "랜덤 가위 바위 보 게임" slide (the second slide in the
# deck) from the final submission - the rest of the deck (the PLC traffic-
# light / "사거리 신호등" slides) is left untouched and simply shifts up by
# one position.
$p.Slides.Item(2).Delete()

$p.Save()
